$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11: TMP_SW_F6_T1 - add Omit / note about re-run
$ws.Range("C11").Value = "Omit"
$ws.Range("D11").Value = "only 1 peak, re-ran 20220708 to get more peaks"

# Row 12: TMP_SW_E3_T1 - add a single space in Action column
$ws.Range("C12").Value = " "

# Row 15: TMP_FW_E3_T1 - add Omit / note about re-run
$ws.Range("C15").Value = "Omit"
$ws.Range("D15").Value = "only 1 peak, re-ran 20220708 to get more peaks"

# Row 37: TMP_FW_POOL_T2 - add Omit / note about re-run
$ws.Range("C37").Value = "Omit"
$ws.Range("D37").Value = "only 1 peak, re-ran 20220708 to get more peaks"

# Row 42: TMP_FW_D5_T2 - add "Omit " (trailing space) / note about re-run
$ws.Range("C42").Value = "Omit "
$ws.Range("D42").Value = "only 1 peak, re-ran 20220708 to get more peaks"

# Switch active sheet to Sheet1 and set the view/selection as in the edit
$ws.Activate()
$ws.Range("D43").Select()
